$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: Architektur weiter ergänzt für bessere GC Optimierung / 1h / 2018-10-31
$ws.Range("B32").Value = "Architektur weiter ergänzt für bessere GC Optimierung"
$ws.Range("D32").Value = "1h"
$ws.Range("F32").Value = 43404
$ws.Range("F32").NumberFormat = "d-mmm"

# Row 33: Funktionalität mit Timern implementiert (Projektilsystem) / 2h / 2018-10-31
$ws.Range("B33").Value = "Funktionalität mit Timern implementiert (Projektilsystem)"
$ws.Range("D33").Value = "2h"
$ws.Range("F33").Value = 43404
$ws.Range("F33").NumberFormat = "d-mmm"

# Row 34: Primitives Lebenssystem ist jetzt auch drinnen / 1h / 2018-10-31
$ws.Range("B34").Value = "Primitives Lebenssystem ist jetzt auch drinnen"
$ws.Range("D34").Value = "1h"
$ws.Range("F34").Value = 43404
$ws.Range("F34").NumberFormat = "d-mmm"

# Update the active selection to follow the newly added last row, like Excel
# would leave it after entering this data.
$ws.Range("B34").Select()
